$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: the rating "3" that used to live in D2 was actually meant to be the
# free-text reason for the /give-reason call -- move it from D (rating) to
# E (reason) and put the real reason text there, clearing D2.
$ws.Range("D2").ClearContents()
$ws.Range("E2").Value = "haha idk"

# Rows 3-6: fill in the missing start_time / end_time / rating / reason data.
$ws.Range("B3").Value = 1689015608877
$ws.Range("C3").Value = 1689015613411
$ws.Range("D3").Value = 3
$ws.Range("E3").Value = "Lepo "

$ws.Range("B4").Value = 1689015617819
$ws.Range("C4").Value = 1689015626587
$ws.Range("D4").Value = 3
$ws.Range("E4").Value = "sdgefc"

$ws.Range("B5").Value = 1689015629330
$ws.Range("C5").Value = 1689015710625
$ws.Range("D5").Value = 3
$ws.Range("E5").Value = "xsefzw"

$ws.Range("B6").Value = 1689015713255
$ws.Range("C6").Value = 1689015730798
$ws.Range("D6").Value = 4
$ws.Range("E6").Value = "sadaws"
